# Generate Report for Handoff
# Update "Latest Handoff Datetime" (column E, row 6 -> file 3afd2715-bbf8-4997-b3c5-ebd3edfae21b)
# on both the zh-cn and de-de localization status sheets to reflect the newest handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-19 00:30:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-19 00:31:00"
